$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Bico")

# Bug fix: the "Obs_relatorio" column (H) was being filled in with a canned
# "Validado com sucesso!..." success message for every row, even though this
# column is meant to hold report text produced elsewhere. Clear out the
# erroneous placeholder text for rows 2-13.
$ws.Range("H2:H13").ClearContents()
